# MM_26 Bill of Parts - "imported MCU to KiCAD"
# Adds a Voltage Regulator line, fills in part details for the IMU,
# STM32 Nucleo, Caster Wheel and OLED Screen rows, and renames the
# generic "STM32 Nucleo Board" part to the specific "STM32 Nucleo L476RG"
# module (moving the old generic name into the Description column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$CURRENCY_FMT = """$""#,##0.00"

# ---------------------------------------------------------------
# Row 5 - IR Sensor Module: give D5 the currency number format even
# though there is still no unit price value for it.
# ---------------------------------------------------------------
$ws.Range("D5").NumberFormat = $CURRENCY_FMT

# ---------------------------------------------------------------
# Row 7 - IMU: part number, description
# ---------------------------------------------------------------
$ws.Range("E7").Value = 4502
$ws.Range("G7").Value = "Gryo + Accel"

# ---------------------------------------------------------------
# Row 8 - STM32 Nucleo board: rename part, add pricing/description
# ---------------------------------------------------------------
$ws.Range("A8").Value = "STM32 Nucleo L476RG"
$ws.Range("D8").NumberFormat = $CURRENCY_FMT
$ws.Range("D8").Value = 14.53
$ws.Range("G8").Value = "STM32 Nucleo Board"
$ws.Range("I8").Value = "https://www.st.com/en/evaluation-tools/nucleo-l476rg.html#st_all-features_sec-nav-tab"

# ---------------------------------------------------------------
# Row 9 - Caster Wheel: pricing, part number, description, url
# ---------------------------------------------------------------
$ws.Range("D9").NumberFormat = $CURRENCY_FMT
$ws.Range("D9").Value = 1.95
$ws.Range("E9").Value = 3948
$ws.Range("G9").Value = "Free spinning wheel"
$ws.Range("I9").Value = "https://www.adafruit.com/product/3948"

# ---------------------------------------------------------------
# Row 10 - OLED Screen: apply currency format, part number, description
# ---------------------------------------------------------------
$ws.Range("D10").NumberFormat = $CURRENCY_FMT
$ws.Range("D10").Value = 17.5
$ws.Range("E10").Value = 661
$ws.Range("G10").Value = "Display"

# ---------------------------------------------------------------
# Row 11 (new) - Voltage Regulator
# ---------------------------------------------------------------
$ws.Range("A11").Value = "Voltage Regulator"
$ws.Range("B11").Value = 1
$ws.Range("D11").NumberFormat = $CURRENCY_FMT
$ws.Range("D11").Value = 0.24
$ws.Range("E11").Value = "AZ1117IH-3.3TRG1DICT-ND"
$ws.Range("G11").Value = "3V3 LDO Regulator"
$ws.Range("I11").Value = "https://www.digikey.com/en/products/detail/diodes-incorporated/AZ1117IH-3-3TRG1/5699672"

# ---------------------------------------------------------------
# Link column (F): hyperlink formula for every new/updated row,
# matching the style already used by F6. Each row's formula points
# at its own URL cell in column I.
# ---------------------------------------------------------------
$ws.Range("F7").Formula = "=HYPERLINK(I7)"
$ws.Range("F8").Formula = "=HYPERLINK(I8)"
$ws.Range("F9").Formula = "=HYPERLINK(I9)"
$ws.Range("F10").Formula = "=HYPERLINK(I10)"
$ws.Range("F11").Formula = "=HYPERLINK(I11)"
$ws.Range("F7:F11").Style = "Hyperlink"

# ---------------------------------------------------------------
# Total Price column (H): extend the total formula down to the new
# STM32 Nucleo row.
# ---------------------------------------------------------------
$ws.Range("H8").Formula = "=(B8+C8)*D8"
$ws.Range("H8").NumberFormat = $CURRENCY_FMT

$wb.Application.Calculate()

# ---------------------------------------------------------------
# Restore the selection state recorded in the saved workbook.
# ---------------------------------------------------------------
[void]$ws.Range("F9").Select()
